$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$col_E = @(2,2,2,2,2,2,3,3,3,3,3,3,3,3,3,3,3,3,3,3,3,3,3,3,3,3,3,3,3,3,2,2,2,2,2,2)
$col_G = @(4.3633995,4.3633995,4.3633995,4.3633995,4.3633995,4.3633995,2.700122,2.700122,2.700122,2.700122,2.700122,2.700122,2.791195666666667,2.791195666666667,2.791195666666667,2.791195666666667,2.791195666666667,2.791195666666667,2.576876,2.576876,2.576876,2.576876,2.576876,2.576876,3.567523333333333,3.567523333333333,3.567523333333333,3.567523333333333,3.567523333333333,3.567523333333333,3.513618,3.513618,3.513618,3.513618,3.513618,3.513618)
$col_H = @(8.726799,8.726799,8.726799,8.726799,8.726799,8.726799,8.100365999999999,8.100365999999999,8.100365999999999,8.100365999999999,8.100365999999999,8.100365999999999,8.373587000000001,8.373587000000001,8.373587000000001,8.373587000000001,8.373587000000001,8.373587000000001,7.730627999999999,7.730627999999999,7.730627999999999,7.730627999999999,7.730627999999999,7.730627999999999,10.70257,10.70257,10.70257,10.70257,10.70257,10.70257,7.027236,7.027236,7.027236,7.027236,7.027236,7.027236)
$col_I = @(0.2236180428734886,0.2236180428734886,0.2236180428734886,0.2236180428734886,0.2236180428734886,0.2236180428734886,0.1383774273154795,0.1383774273154795,0.1383774273154795,0.1383774273154795,0.1383774273154795,0.1383774273154795,0.1430448237107242,0.1430448237107242,0.1430448237107242,0.1430448237107242,0.1430448237107242,0.1430448237107242,0.132061244414513,0.132061244414513,0.132061244414513,0.132061244414513,0.132061244414513,0.132061244414513,0.1828305168265029,0.1828305168265029,0.1828305168265029,0.1828305168265029,0.1828305168265029,0.1828305168265029,0.1800679448592918,0.1800679448592918,0.1800679448592918,0.1800679448592918,0.1800679448592918,0.1800679448592918)
$col_J = @(0.1722580872859944,0.1722580872859944,0.1722580872859944,0.1722580872859944,0.1722580872859944,0.1722580872859944,0.159892940524527,0.159892940524527,0.159892940524527,0.159892940524527,0.159892940524527,0.159892940524527,0.1652860436390099,0.1652860436390099,0.1652860436390099,0.1652860436390099,0.1652860436390099,0.1652860436390099,0.1525946905388278,0.1525946905388278,0.1525946905388278,0.1525946905388278,0.1525946905388278,0.1525946905388278,0.2112577861876348,0.2112577861876348,0.2112577861876348,0.2112577861876348,0.2112577861876348,0.2112577861876348,0.1387104518240058,0.1387104518240058,0.1387104518240058,0.1387104518240058,0.1387104518240058,0.1387104518240058)
$col_K = @(2,3,3,3,3,2,2,3,3,3,3,2,2,3,3,3,3,2,2,3,3,3,3,2,2,3,3,3,3,2,2,3,3,3,3,2)
$col_M = @(2.019046,139.6948166666666,186.3548536666667,143.6051993333333,6.253715333333335,29.9633245,2.019046,139.6948166666666,186.3548536666667,143.6051993333333,6.253715333333335,29.9633245,2.019046,139.6948166666666,186.3548536666667,143.6051993333333,6.253715333333335,29.9633245,2.019046,139.6948166666666,186.3548536666667,143.6051993333333,6.253715333333335,29.9633245,2.019046,139.6948166666666,186.3548536666667,143.6051993333333,6.253715333333335,29.9633245,2.019046,139.6948166666666,186.3548536666667,143.6051993333333,6.253715333333335,29.9633245)
$col_N = @(4.038092,419.0844499999999,559.064561,430.815598,18.761146,59.926649,4.038092,419.0844499999999,559.064561,430.815598,18.761146,59.926649,4.038092,419.0844499999999,559.064561,430.815598,18.761146,59.926649,4.038092,419.0844499999999,559.064561,430.815598,18.761146,59.926649,4.038092,419.0844499999999,559.064561,430.815598,18.761146,59.926649,4.038092,419.0844499999999,559.064561,430.815598,18.761146,59.926649)
$col_O = @(0.003975353327590414,0.2750488370661026,0.3669190239530987,0.2827480934208787,0.01231310631861279,0.05899558591371687,0.003975353327590414,0.2750488370661026,0.3669190239530987,0.2827480934208787,0.01231310631861279,0.05899558591371687,0.003975353327590414,0.2750488370661026,0.3669190239530987,0.2827480934208787,0.01231310631861279,0.05899558591371687,0.003975353327590414,0.2750488370661026,0.3669190239530987,0.2827480934208787,0.01231310631861279,0.05899558591371687,0.003975353327590414,0.2750488370661026,0.3669190239530987,0.2827480934208787,0.01231310631861279,0.05899558591371687,0.003975353327590414,0.2750488370661026,0.3669190239530987,0.2827480934208787,0.01231310631861279,0.05899558591371687)
$col_P = @(0.002707057536954368,0.2809459811695414,0.3747858972750337,0.2888103122968479,0.0125771036621259,0.04017364805949665,0.002707057536954368,0.2809459811695414,0.3747858972750337,0.2888103122968479,0.0125771036621259,0.04017364805949665,0.002707057536954368,0.2809459811695414,0.3747858972750337,0.2888103122968479,0.0125771036621259,0.04017364805949665,0.002707057536954368,0.2809459811695414,0.3747858972750337,0.2888103122968479,0.0125771036621259,0.04017364805949665,0.002707057536954368,0.2809459811695414,0.3747858972750337,0.2888103122968479,0.0125771036621259,0.04017364805949665,0.002707057536954368,0.2809459811695414,0.3747858972750337,0.2888103122968479,0.0125771036621259,0.04017364805949665)
$col_Q = @(8.809904306877,609.5442931959249,813.1406753117066,626.606854968467,27.287458358609,130.7419551416377,5.451670523611999,377.1930477676332,503.1808401921473,387.7515580343186,16.88579435327067,80.90463167558899,5.635552446000667,389.9155669357944,520.1528600167009,400.8302100900029,17.4553431389669,83.63350150332717,5.202831180295999,359.9762203927332,480.2133498971453,370.0527916372827,16.11504895329867,77.211771784262,7.202993716073333,498.3645180040555,664.8252887357523,512.3148994096512,22.31027537169112,106.8948592979883,7.094156368428,490.8342223466999,654.7797682305661,504.573813271188,21.973166762076,105.279676303041)
$col_R = @(35.239617227508,3657.26575917555,4878.844051870239,3759.641129810802,163.724750151654,522.9678205665509,32.710023141672,3394.737429908699,4528.627561729326,3489.764022308868,151.972149179436,485.427790053534,33.813314676004,3509.24010242215,4681.375740150307,3607.471890810026,157.098088250702,501.801009019963,31.216987081776,3239.785983534599,4321.920149074308,3330.475124735544,145.035440579688,463.2706307055719,43.21796229644,4485.280662036499,5983.42759862177,4610.83409468686,200.79247834522,641.3691557879299,28.376625473712,2945.0053340802,3928.678609383396,3027.442879627128,131.839000572456,421.118705212164)
$col_S = @(0.0008889607308463787,0.0615058826393509,0.08204971402944261,0.06322757527698721,0.002753432736661377,0.01319247746020012,0.0005500991661419922,0.03806055045932178,0.05077331056773661,0.03912595374593806,0.001703855974631613,0.008163657401709485,0.0005686537159330118,0.03934431240996036,0.05248586709748201,0.04044565117793298,0.001761326122677371,0.008439013186738516,0.0005249901074289652,0.03632329169771412,0.04845578290260471,0.03734006507299221,0.001626084143044207,0.007791030490728762,0.0007268159034513134,0.05028732103332412,0.06708399478282104,0.05169498005184758,0.002251211591971654,0.01078619346308721,0.0007158337037887527,0.04952747882643128,0.06607035457301173,0.05091386809518066,0.002217195749626564,0.01062321391125278)
$col_T = @(0.0004663125534888947,0.04839521734695221,0.06455990180636249,0.04974991198472575,0.002166507820435484,0.006920235774029594,0.0004328393897527173,0.04492127905774636,0.05992561918242843,0.04617873008694997,0.002010990087819107,0.006423482719830376,0.0004474388301863505,0.04643644970379329,0.06194687817228672,0.0477363138816929,0.002078819704750494,0.006640143346300191,0.0004130826071223534,0.04287086505469354,0.05719033801300069,0.04407092022935974,0.001919199241196861,0.006130285393454673,0.0005718869823395312,0.05935202602019026,0.07917643895266993,0.06101342720399153,0.002657011076313113,0.008486995952130432,0.0003754971740645209,0.03897014398616573,0.05198672114828537,0.04006100891012801,0.001744575731610843,0.005572504873751376)

for ($i = 0; $i -lt 36; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 5).Value = $col_E[$i]
    $ws.Cells.Item($r, 7).Value = $col_G[$i]
    $ws.Cells.Item($r, 8).Value = $col_H[$i]
    $ws.Cells.Item($r, 9).Value = $col_I[$i]
    $ws.Cells.Item($r, 10).Value = $col_J[$i]
    $ws.Cells.Item($r, 11).Value = $col_K[$i]
    $ws.Cells.Item($r, 13).Value = $col_M[$i]
    $ws.Cells.Item($r, 14).Value = $col_N[$i]
    $ws.Cells.Item($r, 15).Value = $col_O[$i]
    $ws.Cells.Item($r, 16).Value = $col_P[$i]
    $ws.Cells.Item($r, 17).Value = $col_Q[$i]
    $ws.Cells.Item($r, 18).Value = $col_R[$i]
    $ws.Cells.Item($r, 19).Value = $col_S[$i]
    $ws.Cells.Item($r, 20).Value = $col_T[$i]
}

Write-Output "applied"